$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-11 (columns A-G)
$data = @(
    @(69587, "Pietro Aragão", "Engenharia", "Outros", 3, 45095, 4153.29),
    @(89328, "Dra. Larissa Moraes", "Recursos Humanos", "Outros", 5, 45105, 6667.69),
    @(42617, "Sra. Maria Luiza Pinto", "Marketing", "Consulta médica", 5, 45097, 5070.68),
    @(73380, "Anthony das Neves", "Financeiro", "Doença", 6, 45106, 2843.46),
    @(63813, "Thales Cunha", "Vendas", "Consulta médica", 5, 45082, 10469.38),
    @(78230, "Elisa da Cruz", "Marketing", "Problemas pessoais", 1, 45079, 6994.97),
    @(29326, "Ana Laura da Rocha", "Jurídico", "Outros", 7, 45098, 9923.51),
    @(49065, "Ana Carolina Melo", "Recursos Humanos", "Doença", 4, 45098, 8246.49),
    @(44194, "Bárbara Monteiro", "Atendimento ao Cliente", "Consulta médica", 6, 45089, 4242.4),
    @(22825, "Sra. Sabrina Silveira", "P&D", "Doença", 8, 45080, 5943.81)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $row++
}
